# Update the 20x5 table of simple arithmetic problems: every cell's
# expression is replaced with a new one (same table shape, only text
# content of the w:t runs changes).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "19+58=77"
$t.Cell(1, 2).Range.Text = "71-7=64"
$t.Cell(1, 3).Range.Text = "76+16=92"
$t.Cell(1, 4).Range.Text = "85-38=47"
$t.Cell(1, 5).Range.Text = "51-7=44"
$t.Cell(2, 1).Range.Text = "40-23=17"
$t.Cell(2, 2).Range.Text = "21-2=19"
$t.Cell(2, 3).Range.Text = "73-68=5"
$t.Cell(2, 4).Range.Text = "56+9=65"
$t.Cell(2, 5).Range.Text = "5+18=23"
$t.Cell(3, 1).Range.Text = "27-18=9"
$t.Cell(3, 2).Range.Text = "37+16=53"
$t.Cell(3, 3).Range.Text = "70-9=61"
$t.Cell(3, 4).Range.Text = "25-16=9"
$t.Cell(3, 5).Range.Text = "42+9=51"
$t.Cell(4, 1).Range.Text = "71-64=7"
$t.Cell(4, 2).Range.Text = "15+76=91"
$t.Cell(4, 3).Range.Text = "42+49=91"
$t.Cell(4, 4).Range.Text = "47+17=64"
$t.Cell(4, 5).Range.Text = "44+39=83"
$t.Cell(5, 1).Range.Text = "47+15=62"
$t.Cell(5, 2).Range.Text = "8+18=26"
$t.Cell(5, 3).Range.Text = "53-46=7"
$t.Cell(5, 4).Range.Text = "57+8=65"
$t.Cell(5, 5).Range.Text = "32-14=18"
$t.Cell(6, 1).Range.Text = "83-15=68"
$t.Cell(6, 2).Range.Text = "15+17=32"
$t.Cell(6, 3).Range.Text = "62-28=34"
$t.Cell(6, 4).Range.Text = "58+35=93"
$t.Cell(6, 5).Range.Text = "82-35=47"
$t.Cell(7, 1).Range.Text = "74-35=39"
$t.Cell(7, 2).Range.Text = "93-19=74"
$t.Cell(7, 3).Range.Text = "95-46=49"
$t.Cell(7, 4).Range.Text = "25+8=33"
$t.Cell(7, 5).Range.Text = "48+26=74"
$t.Cell(8, 1).Range.Text = "60-17=43"
$t.Cell(8, 2).Range.Text = "34+48=82"
$t.Cell(8, 3).Range.Text = "30-19=11"
$t.Cell(8, 4).Range.Text = "98-79=19"
$t.Cell(8, 5).Range.Text = "29+22=51"
$t.Cell(9, 1).Range.Text = "56+18=74"
$t.Cell(9, 2).Range.Text = "22-14=8"
$t.Cell(9, 3).Range.Text = "9+63=72"
$t.Cell(9, 4).Range.Text = "61-55=6"
$t.Cell(9, 5).Range.Text = "66-58=8"
$t.Cell(10, 1).Range.Text = "70-16=54"
$t.Cell(10, 2).Range.Text = "16+18=34"
$t.Cell(10, 3).Range.Text = "86+8=94"
$t.Cell(10, 4).Range.Text = "20-3=17"
$t.Cell(10, 5).Range.Text = "37+28=65"
$t.Cell(11, 1).Range.Text = "93-36=57"
$t.Cell(11, 2).Range.Text = "80-47=33"
$t.Cell(11, 3).Range.Text = "34-29=5"
$t.Cell(11, 4).Range.Text = "75-46=29"
$t.Cell(11, 5).Range.Text = "29+5=34"
$t.Cell(12, 1).Range.Text = "90-18=72"
$t.Cell(12, 2).Range.Text = "53-39=14"
$t.Cell(12, 3).Range.Text = "90-11=79"
$t.Cell(12, 4).Range.Text = "72-59=13"
$t.Cell(12, 5).Range.Text = "59+9=68"
$t.Cell(13, 1).Range.Text = "79+7=86"
$t.Cell(13, 2).Range.Text = "8+58=66"
$t.Cell(13, 3).Range.Text = "20-12=8"
$t.Cell(13, 4).Range.Text = "19+14=33"
$t.Cell(13, 5).Range.Text = "83-37=46"
$t.Cell(14, 1).Range.Text = "64-28=36"
$t.Cell(14, 2).Range.Text = "54+17=71"
$t.Cell(14, 3).Range.Text = "25+6=31"
$t.Cell(14, 4).Range.Text = "72-48=24"
$t.Cell(14, 5).Range.Text = "66-9=57"
$t.Cell(15, 1).Range.Text = "29+28=57"
$t.Cell(15, 2).Range.Text = "89+7=96"
$t.Cell(15, 3).Range.Text = "16+55=71"
$t.Cell(15, 4).Range.Text = "50-37=13"
$t.Cell(15, 5).Range.Text = "80-54=26"
$t.Cell(16, 1).Range.Text = "41-8=33"
$t.Cell(16, 2).Range.Text = "50-26=24"
$t.Cell(16, 3).Range.Text = "62-45=17"
$t.Cell(16, 4).Range.Text = "55-28=27"
$t.Cell(16, 5).Range.Text = "57+34=91"
$t.Cell(17, 1).Range.Text = "35+37=72"
$t.Cell(17, 2).Range.Text = "8+17=25"
$t.Cell(17, 3).Range.Text = "18+63=81"
$t.Cell(17, 4).Range.Text = "44+39=83"
$t.Cell(17, 5).Range.Text = "38+33=71"
$t.Cell(18, 1).Range.Text = "85-56=29"
$t.Cell(18, 2).Range.Text = "53-14=39"
$t.Cell(18, 3).Range.Text = "6+27=33"
$t.Cell(18, 4).Range.Text = "36+27=63"
$t.Cell(18, 5).Range.Text = "38+8=46"
$t.Cell(19, 1).Range.Text = "19+2=21"
$t.Cell(19, 2).Range.Text = "73-58=15"
$t.Cell(19, 3).Range.Text = "6+49=55"
$t.Cell(19, 4).Range.Text = "49+23=72"
$t.Cell(19, 5).Range.Text = "91-25=66"
$t.Cell(20, 1).Range.Text = "75+6=81"
$t.Cell(20, 2).Range.Text = "88+5=93"
$t.Cell(20, 3).Range.Text = "8+66=74"
$t.Cell(20, 4).Range.Text = "54-6=48"
$t.Cell(20, 5).Range.Text = "27+37=64"